$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide-number placeholder preview glyph: "<N>" -> "<#>" (guillemet forms)
#    across every slide layout and the notes master.
# ---------------------------------------------------------------------------
$master = $p.SlideMaster
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.PlaceholderFormat.Type -eq 13) {
            $shp.TextFrame.TextRange.Text = "‹#›"
        }
    }
}

if ($p.HasNotesMaster) {
    $nm = $p.NotesMaster
    for ($i = 1; $i -le $nm.Shapes.Count; $i++) {
        $shp = $nm.Shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.PlaceholderFormat.Type -eq 13) {
            $shp.TextFrame.TextRange.Text = "‹#›"
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 6 ("Requirements and Entities handled by Document DB") text edits.
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$shp6 = $s6.Shapes.Item(2)
$tr6 = $shp6.TextFrame.TextRange

$tr6.Paragraphs(1).Runs(2).Text = "Movie, User, Review, Worker"
$tr6.Paragraphs(11).Runs(1).Text = "Most positive reviewed movies for each genre/year"
$tr6.Paragraphs(12).Runs(1).Text = "Best production houses"

# ---------------------------------------------------------------------------
# 3) Slide 7 ("Requirements and Entities handled by Graph DB") text edit.
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$shp7 = $s7.Shapes.Item(2)
$tr7 = $shp7.TextFrame.TextRange

$tr7.Paragraphs(9).Runs(1).Text = "Last reviewed movies of followed top critics"
